# Apply updated profit/price figures per scheduled-runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2118.606
$ws.Range("I15").Value = 2118.606
$ws.Range("K15").Value = 6355.818000000001
$ws.Range("M15").Value = -6186.818000000001

$ws.Range("H94").Value = 9334.166999999999
$ws.Range("I94").Value = 3230
$ws.Range("K94").Value = 3230
$ws.Range("M94").Value = -2779

$ws.Range("H100").Value = 29413460
$ws.Range("I100").Value = 1363.75
$ws.Range("J100").Value = 55557544
$ws.Range("K100").Value = 1363.75
$ws.Range("L100").Value = 55557544
$ws.Range("M100").Value = -822.75
$ws.Range("N100").Value = -55558626

$ws.Range("H113").Value = 41669588
$ws.Range("I113").Value = 250000940
$ws.Range("J113").Value = 3315
$ws.Range("K113").Value = 250000940
$ws.Range("L113").Value = 3315
$ws.Range("M113").Value = -249997686
$ws.Range("N113").Value = -9823

$ws.Range("H141").Value = 1431.9436
$ws.Range("I141").Value = 950.1778
$ws.Range("J141").Value = 2265.7693
$ws.Range("K141").Value = 2850.5334
$ws.Range("L141").Value = 6797.3079
$ws.Range("M141").Value = 2329.4666
$ws.Range("N141").Value = -17157.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2648782.5
$ws.Range("I32").Value = 4277.6626
$ws.Range("J32").Value = 16752809
$ws.Range("K32").Value = 4277.6626
$ws.Range("L32").Value = 16752809
$ws.Range("M32").Value = -3990.6626
$ws.Range("N32").Value = -16753383

$ws.Range("H45").Value = 4466.3335
$ws.Range("I45").Value = 3250.1
$ws.Range("J45").Value = 5986.625
$ws.Range("K45").Value = 3250.1
$ws.Range("L45").Value = 5986.625
$ws.Range("M45").Value = -2873.1
$ws.Range("N45").Value = -6740.625

$ws.Range("H61").Value = 1236.4565
$ws.Range("I61").Value = 1376.4193
$ws.Range("J61").Value = 947.2
$ws.Range("K61").Value = 1376.4193
$ws.Range("L61").Value = 947.2
$ws.Range("M61").Value = -1164.4193
$ws.Range("N61").Value = -1371.2

$ws.Range("H132").Value = 56988.99
$ws.Range("I132").Value = 67631.56
$ws.Range("J132").Value = 3776.1333
$ws.Range("K132").Value = 202894.68
$ws.Range("L132").Value = 11328.3999
$ws.Range("M132").Value = -200364.68
$ws.Range("N132").Value = -16388.3999

$ws.Range("H136").Value = 1236.4565
$ws.Range("I136").Value = 1376.4193
$ws.Range("J136").Value = 947.2
$ws.Range("K136").Value = 4129.257900000001
$ws.Range("L136").Value = 2841.6
$ws.Range("M136").Value = -1579.257900000001
$ws.Range("N136").Value = -7941.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1720.125
$ws.Range("I31").Value = 1346.6216
$ws.Range("J31").Value = 6326.6665
$ws.Range("K31").Value = 1346.6216
$ws.Range("L31").Value = 6326.6665
$ws.Range("M31").Value = -1051.6216
$ws.Range("N31").Value = -6916.6665

$ws.Range("H34").Value = 1720.125
$ws.Range("I34").Value = 1346.6216
$ws.Range("J34").Value = 6326.6665
$ws.Range("K34").Value = 1346.6216
$ws.Range("L34").Value = 6326.6665
$ws.Range("M34").Value = -1144.6216
$ws.Range("N34").Value = -6730.6665

$ws.Range("H107").Value = 1140.7838
$ws.Range("I107").Value = 1146.4073
$ws.Range("J107").Value = 1125.6
$ws.Range("K107").Value = 1146.4073
$ws.Range("L107").Value = 1125.6
$ws.Range("M107").Value = 773.5926999999999
$ws.Range("N107").Value = -4965.6

$ws.Range("H132").Value = 1817.6349
$ws.Range("I132").Value = 1633.7115
$ws.Range("J132").Value = 2687.0908
$ws.Range("K132").Value = 4901.1345
$ws.Range("L132").Value = 8061.2724
$ws.Range("M132").Value = -2371.1345
$ws.Range("N132").Value = -13121.2724

$ws.Range("H134").Value = 3647.5098
$ws.Range("I134").Value = 3708.6738
$ws.Range("J134").Value = 3084.8
$ws.Range("K134").Value = 11126.0214
$ws.Range("L134").Value = 9254.400000000001
$ws.Range("M134").Value = -8591.0214
$ws.Range("N134").Value = -14324.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 360627.3
$ws.Range("I5").Value = 274.3889
$ws.Range("K5").Value = 823.1667
$ws.Range("M5").Value = -711.1667

$ws.Range("H23").Value = 22.8
$ws.Range("J23").Value = 21.333334
$ws.Range("L23").Value = 64.00000199999999
$ws.Range("N23").Value = -534.000002

$ws.Range("H70").Value = 6439.7915
$ws.Range("I70").Value = 5737
$ws.Range("K70").Value = 17211
$ws.Range("M70").Value = -16896

$ws.Range("H73").Value = 6439.7915
$ws.Range("I73").Value = 5737
$ws.Range("K73").Value = 17211
$ws.Range("M73").Value = -16119

$ws.Range("H122").Value = 43344.9
$ws.Range("I122").Value = 311.84616
$ws.Range("J122").Value = 49849.895
$ws.Range("K122").Value = 2806.61544
$ws.Range("L122").Value = 448649.055
$ws.Range("M122").Value = -356.61544
$ws.Range("N122").Value = -453549.055

$ws.Range("H135").Value = 360627.3
$ws.Range("I135").Value = 274.3889
$ws.Range("K135").Value = 2469.5001
$ws.Range("M135").Value = 65.49990000000025

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4985
$ws.Range("I126").Value = 2884.6155
$ws.Range("J126").Value = 6591.1763
$ws.Range("K126").Value = 8653.8465
$ws.Range("L126").Value = 19773.5289
$ws.Range("M126").Value = -6183.8465
$ws.Range("N126").Value = -24713.5289

$ws.Range("H132").Value = 1872.0635
$ws.Range("I132").Value = 1540.7
$ws.Range("J132").Value = 2448.348
$ws.Range("K132").Value = 4622.1
$ws.Range("L132").Value = 7345.044
$ws.Range("M132").Value = -2092.1
$ws.Range("N132").Value = -12405.044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4764440.5
$ws.Range("I7").Value = 12501963
$ws.Range("K7").Value = 12501963
$ws.Range("M7").Value = -12501851

$ws.Range("H40").Value = 1868.5
$ws.Range("I40").Value = 1550
$ws.Range("K40").Value = 1550
$ws.Range("M40").Value = -1414

$ws.Range("H46").Value = 1332.174
$ws.Range("I46").Value = 1141.1111
$ws.Range("J46").Value = 2020
$ws.Range("K46").Value = 1141.1111
$ws.Range("L46").Value = 2020
$ws.Range("M46").Value = -953.1111000000001
$ws.Range("N46").Value = -2396

$ws.Range("H126").Value = 4764440.5
$ws.Range("I126").Value = 12501963
$ws.Range("K126").Value = 37505889
$ws.Range("M126").Value = -37503419

$ws.Range("H136").Value = 1297.8591
$ws.Range("I136").Value = 1099.9824
$ws.Range("K136").Value = 3299.947200000001
$ws.Range("M136").Value = -749.9472000000005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1909.8572
$ws.Range("I132").Value = 1839.9656
$ws.Range("J132").Value = 2720.6
$ws.Range("K132").Value = 5519.8968
$ws.Range("L132").Value = 8161.799999999999
$ws.Range("M132").Value = -2989.8968
$ws.Range("N132").Value = -13221.8
